$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Row = 2; Col = 4; Value = "256.44" },
    @{ Row = 2; Col = 5; Value = "0.65%" },
    @{ Row = 3; Col = 4; Value = "26.87" },
    @{ Row = 3; Col = 5; Value = "-3.82%" },
    @{ Row = 4; Col = 4; Value = "4.744" },
    @{ Row = 4; Col = 5; Value = "-9.23%" },
    @{ Row = 5; Col = 4; Value = "0.05927" },
    @{ Row = 5; Col = 5; Value = "1.02%" },
    @{ Row = 6; Col = 4; Value = "6.660" },
    @{ Row = 6; Col = 5; Value = "-0.83%" },
    @{ Row = 7; Col = 4; Value = "0.8670" },
    @{ Row = 7; Col = 5; Value = "0.38%" },
    @{ Row = 8; Col = 4; Value = "0.9381" },
    @{ Row = 8; Col = 5; Value = "-3.25%" },
    @{ Row = 9; Col = 5; Value = "-0.64%" },
    @{ Row = 10; Col = 4; Value = "0.03794" },
    @{ Row = 10; Col = 5; Value = "9.09%" },
    @{ Row = 11; Col = 4; Value = "0.07113" },
    @{ Row = 11; Col = 5; Value = "-0.85%" },
    @{ Row = 12; Col = 4; Value = "0.03164" },
    @{ Row = 12; Col = 5; Value = "-0.27%" },
    @{ Row = 13; Col = 4; Value = "0.09258" },
    @{ Row = 13; Col = 5; Value = "0.40%" },
    @{ Row = 14; Col = 4; Value = "0.001539" },
    @{ Row = 14; Col = 5; Value = "0.09%" },
    @{ Row = 15; Col = 4; Value = "0.0006068" },
    @{ Row = 15; Col = 5; Value = "-0.38%" },
    @{ Row = 16; Col = 4; Value = "0.006082" },
    @{ Row = 16; Col = 5; Value = "4.90%" },
    @{ Row = 17; Col = 4; Value = "3.497" },
    @{ Row = 17; Col = 5; Value = "-0.08%" },
    @{ Row = 18; Col = 5; Value = "-0.53%" },
    @{ Row = 19; Col = 5; Value = "-0.14%" },
    @{ Row = 20; Col = 4; Value = "0.3146" },
    @{ Row = 20; Col = 5; Value = "-1.09%" },
    @{ Row = 21; Col = 5; Value = "0.37%" },
    @{ Row = 22; Col = 4; Value = "3.820" },
    @{ Row = 22; Col = 5; Value = "7.42%" },
    @{ Row = 23; Col = 4; Value = "0.04223" },
    @{ Row = 23; Col = 5; Value = "1.38%" },
    @{ Row = 24; Col = 5; Value = "-0.02%" },
    @{ Row = 25; Col = 4; Value = "0.001224" },
    @{ Row = 25; Col = 5; Value = "-0.26%" },
    @{ Row = 26; Col = 4; Value = "0.004286" },
    @{ Row = 27; Col = 4; Value = "0.0001200" },
    @{ Row = 27; Col = 5; Value = "0.00%" },
    @{ Row = 28; Col = 4; Value = "0.0001493" },
    @{ Row = 28; Col = 5; Value = "1.79%" },
    @{ Row = 40; Col = 4; Value = "0.03826" },
    @{ Row = 40; Col = 5; Value = "0.36%" },
    @{ Row = 41; Col = 4; Value = "0.006178" },
    @{ Row = 41; Col = 5; Value = "61.66%" },
    @{ Row = 42; Col = 5; Value = "-0.06%" },
    @{ Row = 43; Col = 4; Value = "0.002251" },
    @{ Row = 43; Col = 5; Value = "-3.94%" },
    @{ Row = 44; Col = 4; Value = "0.01115" },
    @{ Row = 44; Col = 5; Value = "14.85%" },
    @{ Row = 45; Col = 4; Value = "0.00005499" },
    @{ Row = 45; Col = 5; Value = "4.93%" },
    @{ Row = 46; Col = 4; Value = "0.00000000750" },
    @{ Row = 47; Col = 4; Value = "0.08850" },
    @{ Row = 47; Col = 5; Value = "-11.53%" },
    @{ Row = 48; Col = 4; Value = "0.002436" },
    @{ Row = 48; Col = 5; Value = "14.16%" },
    @{ Row = 49; Col = 4; Value = "0.00002099" },
)

foreach ($u in $updates) {
    $cell = $ws.Cells.Item($u.Row, $u.Col)
    $cell.NumberFormat = "@"
    $cell.Value = $u.Value
    $cell.Style = "Normal"
}
